# Re-number the Word-generated "_Toc" bookmarks that mark the TOC entries.
# (Word mints a fresh set of these hidden bookmark names whenever the TOC
# field is refreshed; the "_Toc244461xx" set is left untouched, only the
# "_Toc247313xx" set present in this document is renumbered.)
$d = $word.ActiveDocument

$renames = @(
    @{ Old = "_Toc24731300"; New = "_Toc25148399" },
    @{ Old = "_Toc24731301"; New = "_Toc25148400" },
    @{ Old = "_Toc24731302"; New = "_Toc25148401" },
    @{ Old = "_Toc24731303"; New = "_Toc25148402" },
    @{ Old = "_Toc24731304"; New = "_Toc25148403" },
    @{ Old = "_Toc24731305"; New = "_Toc25148404" },
    @{ Old = "_Toc24731306"; New = "_Toc25148405" },
    @{ Old = "_Toc24731307"; New = "_Toc25148406" },
    @{ Old = "_Toc24731308"; New = "_Toc25148407" },
    @{ Old = "_Toc24731309"; New = "_Toc25148408" },
    @{ Old = "_Toc24731310"; New = "_Toc25148409" },
    @{ Old = "_Toc24731311"; New = "_Toc25148410" },
    @{ Old = "_Toc24731312"; New = "_Toc25148411" },
    @{ Old = "_Toc24731313"; New = "_Toc25148412" },
    @{ Old = "_Toc24731314"; New = "_Toc25148413" },
    @{ Old = "_Toc24731315"; New = "_Toc25148414" },
    @{ Old = "_Toc24731316"; New = "_Toc25148415" },
    @{ Old = "_Toc24731317"; New = "_Toc25148416" },
    @{ Old = "_Toc24731318"; New = "_Toc25148417" },
    @{ Old = "_Toc24731319"; New = "_Toc25148418" },
    @{ Old = "_Toc24731320"; New = "_Toc25148419" },
    @{ Old = "_Toc24731321"; New = "_Toc25148420" },
    @{ Old = "_Toc24731322"; New = "_Toc25148421" },
    @{ Old = "_Toc24731323"; New = "_Toc25148422" }
)

foreach ($pair in $renames) {
    $bm = $d.Bookmarks.Item($pair.Old)
    $start = $bm.Start
    $end = $bm.End
    $bm.Delete()
    $r = $d.Range($start, $end)
    $d.Bookmarks.Add($pair.New, $r) | Out-Null
}

# The footer page-number field in section 1's primary footer was
# recalculated (it now reports page 4 instead of page 1); update the
# cached field result text to match.
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Select()
$word.Selection.Text = "4"

Write-Output "done"
